$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1524
$ws.Range("I38").Value = 182
$ws.Range("J38").Value = 2500
$ws.Range("K38").Value = 546
$ws.Range("L38").Value = 7500
$ws.Range("M38").Value = -174
$ws.Range("N38").Value = -8244

$ws.Range("H116").Value = 48751.914
$ws.Range("I116").Value = 64887.707
$ws.Range("J116").Value = 3033.8333
$ws.Range("K116").Value = 64887.707
$ws.Range("L116").Value = 3033.8333
$ws.Range("M116").Value = -61445.707
$ws.Range("N116").Value = -9917.8333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 310765.5
$ws.Range("I61").Value = 245859.83
$ws.Range("J61").Value = 437754.88
$ws.Range("K61").Value = 245859.83
$ws.Range("L61").Value = 437754.88
$ws.Range("M61").Value = -245647.83
$ws.Range("N61").Value = -438178.88

$ws.Range("H63").Value = 2806
$ws.Range("I63").Value = 2597.0476
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 2597.0476
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -1911.0476
$ws.Range("N63").Value = -6372

$ws.Range("H66").Value = 2806
$ws.Range("I66").Value = 2597.0476
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 12985.238
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -9553.237999999999
$ws.Range("N66").Value = -31864

$ws.Range("H74").Value = 146266.27
$ws.Range("I74").Value = 186074.8
$ws.Range("J74").Value = 48554.41
$ws.Range("K74").Value = 186074.8
$ws.Range("L74").Value = 48554.41
$ws.Range("M74").Value = -185200.8
$ws.Range("N74").Value = -50302.41

$ws.Range("H77").Value = 146266.27
$ws.Range("I77").Value = 186074.8
$ws.Range("J77").Value = 48554.41
$ws.Range("K77").Value = 930374
$ws.Range("L77").Value = 242772.05
$ws.Range("M77").Value = -926006
$ws.Range("N77").Value = -251508.05

$ws.Range("H122").Value = 4030.0952
$ws.Range("I122").Value = 4191.7856
$ws.Range("J122").Value = 3706.7144
$ws.Range("K122").Value = 12575.3568
$ws.Range("L122").Value = 11120.1432
$ws.Range("M122").Value = -10125.3568
$ws.Range("N122").Value = -16020.1432

$ws.Range("H136").Value = 310765.5
$ws.Range("I136").Value = 245859.83
$ws.Range("J136").Value = 437754.88
$ws.Range("K136").Value = 737579.49
$ws.Range("L136").Value = 1313264.64
$ws.Range("M136").Value = -735029.49
$ws.Range("N136").Value = -1318364.64

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2419.8135
$ws.Range("I134").Value = 2049.3333
$ws.Range("J134").Value = 4036.4546
$ws.Range("K134").Value = 6147.999899999999
$ws.Range("L134").Value = 12109.3638
$ws.Range("M134").Value = -3612.999899999999
$ws.Range("N134").Value = -17179.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2753.2
$ws.Range("I31").Value = 2015.65
$ws.Range("M31").Value = -1720.65

$ws.Range("H34").Value = 2753.2
$ws.Range("I34").Value = 2015.65
$ws.Range("K34").Value = 2015.65
$ws.Range("M34").Value = -1813.65

$ws.Range("H62").Value = 6175936
$ws.Range("I62").Value = 13890962
$ws.Range("J62").Value = 3915.4
$ws.Range("K62").Value = 13890962
$ws.Range("L62").Value = 3915.4
$ws.Range("M62").Value = -13890338
$ws.Range("N62").Value = -5163.4

$ws.Range("H65").Value = 6175936
$ws.Range("I65").Value = 13890962
$ws.Range("J65").Value = 3915.4
$ws.Range("K65").Value = 69454810
$ws.Range("L65").Value = 19577
$ws.Range("M65").Value = -69451690
$ws.Range("N65").Value = -25817

$ws.Range("H105").Value = 887.38464
$ws.Range("I105").Value = 705.36365
$ws.Range("K105").Value = 705.36365
$ws.Range("M105").Value = 1041.63635

$ws.Range("H122").Value = 1982.909
$ws.Range("I122").Value = 853
$ws.Range("J122").Value = 2628.5715
$ws.Range("K122").Value = 2559
$ws.Range("L122").Value = 7885.7145
$ws.Range("M122").Value = -109
$ws.Range("N122").Value = -12785.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 15.318182
$ws.Range("I12").Value = 22.11111
$ws.Range("J12").Value = 10.615385
$ws.Range("K12").Value = 66.33333
$ws.Range("L12").Value = 31.846155
$ws.Range("M12").Value = 106.66667
$ws.Range("N12").Value = -377.846155

$ws.Range("H14").Value = 460.26923
$ws.Range("I14").Value = 460.26923
$ws.Range("K14").Value = 1380.80769
$ws.Range("M14").Value = -1207.80769

$ws.Range("H98").Value = 4425.1875
$ws.Range("I98").Value = 230.4
$ws.Range("J98").Value = 11416.5
$ws.Range("K98").Value = 691.2
$ws.Range("L98").Value = 34249.5
$ws.Range("M98").Value = 806.8
$ws.Range("N98").Value = -37245.5

$ws.Range("H122").Value = 21740032
$ws.Range("I122").Value = 41667004
$ws.Range("J122").Value = 1515.8182
$ws.Range("K122").Value = 375003036
$ws.Range("L122").Value = 13642.3638
$ws.Range("M122").Value = -375000586
$ws.Range("N122").Value = -18542.3638

$ws.Range("H132").Value = 6458.3076
$ws.Range("I132").Value = 1687.6
$ws.Range("J132").Value = 9440
$ws.Range("K132").Value = 15188.4
$ws.Range("L132").Value = 84960
$ws.Range("M132").Value = -12658.4
$ws.Range("N132").Value = -90020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5700
$ws.Range("I80").Value = 6781.25
$ws.Range("J80").Value = 2816.6667
$ws.Range("K80").Value = 6781.25
$ws.Range("L80").Value = 2816.6667
$ws.Range("M80").Value = -5783.25
$ws.Range("N80").Value = -4812.6667

$ws.Range("H83").Value = 5700
$ws.Range("I83").Value = 6781.25
$ws.Range("J83").Value = 2816.6667
$ws.Range("K83").Value = 33906.25
$ws.Range("L83").Value = 14083.3335
$ws.Range("M83").Value = -28914.25
$ws.Range("N83").Value = -24067.3335

$ws.Range("H102").Value = 4826.476
$ws.Range("I102").Value = 2745.0588
$ws.Range("K102").Value = 2745.0588
$ws.Range("M102").Value = -1123.0588

$ws.Range("H122").Value = 216221.64
$ws.Range("I122").Value = 601484.6
$ws.Range("J122").Value = 2186.6667
$ws.Range("K122").Value = 1804453.8
$ws.Range("L122").Value = 6560.000100000001
$ws.Range("M122").Value = -1802003.8
$ws.Range("N122").Value = -11460.0001

$ws.Range("H126").Value = 1961
$ws.Range("I126").Value = 1639.8182
$ws.Range("J126").Value = 2667.6
$ws.Range("K126").Value = 4919.4546
$ws.Range("L126").Value = 8002.799999999999
$ws.Range("M126").Value = -2449.4546
$ws.Range("N126").Value = -12942.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2147.5186
$ws.Range("I122").Value = 2040.1765
$ws.Range("J122").Value = 2330
$ws.Range("K122").Value = 6120.529500000001
$ws.Range("L122").Value = 6990
$ws.Range("M122").Value = -3670.529500000001
$ws.Range("N122").Value = -11890

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 997.6667
$ws.Range("I122").Value = 914
$ws.Range("J122").Value = 1165
$ws.Range("K122").Value = 2742
$ws.Range("L122").Value = 3495
$ws.Range("M122").Value = -292
$ws.Range("N122").Value = -8395
